# Anil changes after payment module testing.
# Clear out the placeholder "Hi this is description / Hi this is conditions /
# Hi, My reason isd that" test data that had been entered in the StageDetails
# template rows while testing the payment module, along with the related
# date/amount cells that were filled in alongside them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that should end up fully empty (cell removed entirely): stage
# description (D), conditions (J), reasons for slippage (L) and the stray
# numeric amount (P).
$ws.Range("D2:D7").ClearContents() | Out-Null
$ws.Range("J2:J7").ClearContents() | Out-Null
$ws.Range("L2:L7").ClearContents() | Out-Null
$ws.Range("P2:P7").ClearContents() | Out-Null

# Date columns that should remain present (keeping their date number format)
# but with no value: stage start/completion dates (E/F) and revised date of
# payment (K).
$ws.Range("E2:F7").ClearContents() | Out-Null
$ws.Range("K2:K7").ClearContents() | Out-Null

# Move the selection to reflect where the author left off reviewing the
# sheet after the cleanup.
$ws.Range("N9").Select() | Out-Null
